$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label for the new column
$ws.Range("N2").Value = "WingLoading"

# Wing loading formula (MTOW / Wing Area) for each data row, rows 3 through 39
for ($r = 3; $r -le 39; $r++) {
    $ws.Range("N$r").Formula = "=C$r/H$r"
}

# Update the active selection to match the saved workbook state
$ws.Range("B3").Select()

$wb.Save()
